$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 501.05264
$ws.Range("I92").Value = 432.85715
$ws.Range("J92").Value = 692
$ws.Range("K92").Value = 432.85715
$ws.Range("L92").Value = 692
$ws.Range("M92").Value = 815.14285
$ws.Range("N92").Value = -3188

$ws.Range("H121").Value = 4802.3447
$ws.Range("J121").Value = 4802.3447
$ws.Range("L121").Value = 14407.0341
$ws.Range("N121").Value = -17901.0341

$ws.Range("H125").Value = 1303.9166
$ws.Range("I125").Value = 715.5
$ws.Range("J125").Value = 1421.6
$ws.Range("K125").Value = 6439.5
$ws.Range("L125").Value = 12794.4
$ws.Range("M125").Value = -3979.5
$ws.Range("N125").Value = -17714.4

$ws.Range("H132").Value = 34067.312
$ws.Range("I132").Value = 40087.332
$ws.Range("J132").Value = 1559.2
$ws.Range("K132").Value = 120261.996
$ws.Range("L132").Value = 4677.6
$ws.Range("M132").Value = -117731.996
$ws.Range("N132").Value = -9737.6

$ws.Range("H137").Value = 19409.475
$ws.Range("I137").Value = 1903.1389
$ws.Range("J137").Value = 49420.332
$ws.Range("K137").Value = 5709.4167
$ws.Range("L137").Value = 148260.996
$ws.Range("M137").Value = -3159.4167
$ws.Range("N137").Value = -153360.996

$ws.Range("H138").Value = 3135.9
$ws.Range("J138").Value = 3684.7222
$ws.Range("L138").Value = 11054.1666
$ws.Range("N138").Value = -21334.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20243.396
$ws.Range("I32").Value = 22471.49
$ws.Range("K32").Value = 22471.49
$ws.Range("M32").Value = -22184.49

$ws.Range("H63").Value = 3127450
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 6252400
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 6252400
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -6253772

$ws.Range("H66").Value = 3127450
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 6252400
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 31262000
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -31268864

$ws.Range("H74").Value = 1069.4138
$ws.Range("I74").Value = 807.0625
$ws.Range("J74").Value = 1392.3077
$ws.Range("K74").Value = 807.0625
$ws.Range("L74").Value = 1392.3077
$ws.Range("M74").Value = 66.9375
$ws.Range("N74").Value = -3140.3077

$ws.Range("H77").Value = 1069.4138
$ws.Range("I77").Value = 807.0625
$ws.Range("J77").Value = 1392.3077
$ws.Range("K77").Value = 4035.3125
$ws.Range("L77").Value = 6961.538500000001
$ws.Range("M77").Value = 332.6875
$ws.Range("N77").Value = -15697.5385

$ws.Range("H122").Value = 2135.2432
$ws.Range("I122").Value = 1938.4482
$ws.Range("J122").Value = 2848.625
$ws.Range("K122").Value = 5815.3446
$ws.Range("L122").Value = 8545.875
$ws.Range("M122").Value = -3365.3446
$ws.Range("N122").Value = -13445.875

$ws.Range("H132").Value = 21948.4
$ws.Range("I132").Value = 1380.0714
$ws.Range("J132").Value = 48126.273
$ws.Range("K132").Value = 4140.2142
$ws.Range("L132").Value = 144378.819
$ws.Range("M132").Value = -1610.2142
$ws.Range("N132").Value = -149438.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 980
$ws.Range("I20").Value = 950
$ws.Range("J20").Value = 1100
$ws.Range("K20").Value = 950
$ws.Range("L20").Value = 1100
$ws.Range("M20").Value = -703
$ws.Range("N20").Value = -1594

$ws.Range("H80").Value = 996.5599999999999
$ws.Range("I80").Value = 1707.125
$ws.Range("J80").Value = 662.17645
$ws.Range("K80").Value = 1707.125
$ws.Range("L80").Value = 662.17645
$ws.Range("M80").Value = -709.125
$ws.Range("N80").Value = -2658.17645

$ws.Range("H83").Value = 996.5599999999999
$ws.Range("I83").Value = 1707.125
$ws.Range("J83").Value = 662.17645
$ws.Range("K83").Value = 8535.625
$ws.Range("L83").Value = 3310.88225
$ws.Range("M83").Value = -3543.625
$ws.Range("N83").Value = -13294.88225

$ws.Range("H134").Value = 25535.535
$ws.Range("I134").Value = 30235.723
$ws.Range("J134").Value = 1363.1428
$ws.Range("K134").Value = 90707.16900000001
$ws.Range("L134").Value = 4089.4284
$ws.Range("M134").Value = -88172.16900000001
$ws.Range("N134").Value = -9159.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10112.673
$ws.Range("I31").Value = 22315.895
$ws.Range("J31").Value = 3086.5757
$ws.Range("K31").Value = 22315.895
$ws.Range("L31").Value = 3086.5757
$ws.Range("M31").Value = -22020.895
$ws.Range("N31").Value = -3676.5757

$ws.Range("H34").Value = 10112.673
$ws.Range("I34").Value = 22315.895
$ws.Range("J34").Value = 3086.5757
$ws.Range("K34").Value = 22315.895
$ws.Range("L34").Value = 3086.5757
$ws.Range("M34").Value = -22113.895
$ws.Range("N34").Value = -3490.5757

$ws.Range("H58").Value = 11122.429
$ws.Range("I58").Value = 944
$ws.Range("J58").Value = 63286.875
$ws.Range("K58").Value = 944
$ws.Range("L58").Value = 63286.875
$ws.Range("M58").Value = -741
$ws.Range("N58").Value = -63692.875

$ws.Range("H62").Value = 5433.3335
$ws.Range("I62").Value = 4225
$ws.Range("K62").Value = 4225
$ws.Range("M62").Value = -3601

$ws.Range("H65").Value = 5433.3335
$ws.Range("I65").Value = 4225
$ws.Range("K65").Value = 21125
$ws.Range("M65").Value = -18005

$ws.Range("H134").Value = 861.75
$ws.Range("I134").Value = 843.0345
$ws.Range("J134").Value = 1042.6666
$ws.Range("K134").Value = 2529.1035
$ws.Range("L134").Value = 3127.9998
$ws.Range("M134").Value = 5.896499999999833
$ws.Range("N134").Value = -8197.9998

$ws.Range("H135").Value = 51637.5
$ws.Range("J135").Value = 51637.5
$ws.Range("L135").Value = 51637.5
$ws.Range("N135").Value = -61777.5

$ws.Range("H136").Value = 11122.429
$ws.Range("I136").Value = 944
$ws.Range("J136").Value = 63286.875
$ws.Range("K136").Value = 2832
$ws.Range("L136").Value = 189860.625
$ws.Range("M136").Value = -282
$ws.Range("N136").Value = -194960.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4763.3213
$ws.Range("I68").Value = 640.7692
$ws.Range("J68").Value = 8336.200000000001
$ws.Range("K68").Value = 1922.3076
$ws.Range("L68").Value = 25008.6
$ws.Range("M68").Value = -1111.3076
$ws.Range("N68").Value = -26630.6

$ws.Range("H71").Value = 4763.3213
$ws.Range("I71").Value = 640.7692
$ws.Range("J71").Value = 8336.200000000001
$ws.Range("K71").Value = 5766.922799999999
$ws.Range("L71").Value = 75025.8
$ws.Range("M71").Value = -1710.922799999999
$ws.Range("N71").Value = -83137.8

$ws.Range("H75").Value = 4915
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4915
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 14745
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -16741

$ws.Range("H78").Value = 4915
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4915
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 44235
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -54219

$ws.Range("H106").Value = 5254.5
$ws.Range("J106").Value = 5254.5
$ws.Range("L106").Value = 15763.5
$ws.Range("N106").Value = -17655.5

$ws.Range("H107").Value = 4681.32
$ws.Range("J107").Value = 768.64703
$ws.Range("L107").Value = 2305.94109
$ws.Range("N107").Value = -6145.94109

$ws.Range("H113").Value = 5966.4443
$ws.Range("I113").Value = 12871
$ws.Range("J113").Value = 442.8
$ws.Range("K113").Value = 38613
$ws.Range("L113").Value = 1328.4
$ws.Range("M113").Value = -36443
$ws.Range("N113").Value = -5668.4

$ws.Range("H131").Value = 108339.03
$ws.Range("I131").Value = 753.3333
$ws.Range("J131").Value = 115758.734
$ws.Range("K131").Value = 2259.9999
$ws.Range("L131").Value = 347276.202
$ws.Range("M131").Value = 2780.0001
$ws.Range("N131").Value = -357356.202

$ws.Range("H133").Value = 2820
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1305.6072
$ws.Range("I102").Value = 1235.8334
$ws.Range("K102").Value = 1235.8334
$ws.Range("M102").Value = 386.1666

$ws.Range("H126").Value = 4055.5366
$ws.Range("I126").Value = 3484.261
$ws.Range("J126").Value = 4785.5
$ws.Range("K126").Value = 10452.783
$ws.Range("L126").Value = 14356.5
$ws.Range("M126").Value = -7982.782999999999
$ws.Range("N126").Value = -19296.5

$ws.Range("H132").Value = 48880.758
$ws.Range("I132").Value = 47168.13
$ws.Range("J132").Value = 52819.8
$ws.Range("K132").Value = 141504.39
$ws.Range("L132").Value = 158459.4
$ws.Range("M132").Value = -138974.39
$ws.Range("N132").Value = -163519.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 72189.44
$ws.Range("I40").Value = 160991.58
$ws.Range("J40").Value = 3121.111
$ws.Range("K40").Value = 160991.58
$ws.Range("L40").Value = 3121.111
$ws.Range("M40").Value = -160855.58
$ws.Range("N40").Value = -3393.111

$ws.Range("H46").Value = 1888.9412
$ws.Range("I46").Value = 1808
$ws.Range("J46").Value = 2266.6667
$ws.Range("K46").Value = 1808
$ws.Range("L46").Value = 2266.6667
$ws.Range("M46").Value = -1620
$ws.Range("N46").Value = -2642.6667

$ws.Range("H93").Value = 2775.6365
$ws.Range("I93").Value = 3079
$ws.Range("K93").Value = 3079
$ws.Range("M93").Value = -1831

$ws.Range("H132").Value = 1801.2667
$ws.Range("J132").Value = 2352.8
$ws.Range("L132").Value = 7058.400000000001
$ws.Range("N132").Value = -12118.4

$ws.Range("H136").Value = 14868.923
$ws.Range("I136").Value = 23886.137
$ws.Range("J136").Value = 3199.5881
$ws.Range("K136").Value = 71658.41099999999
$ws.Range("L136").Value = 9598.764299999999
$ws.Range("M136").Value = -69108.41099999999
$ws.Range("N136").Value = -14698.7643

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 21514.5
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H45").Value = 20812.5
$ws.Range("J45").Value = 21000
$ws.Range("L45").Value = 21000
$ws.Range("N45").Value = -21982

$ws.Range("H122").Value = 1481.0714
$ws.Range("I122").Value = 1362
$ws.Range("J122").Value = 2028.8
$ws.Range("K122").Value = 4086
$ws.Range("L122").Value = 6086.4
$ws.Range("M122").Value = -1636
$ws.Range("N122").Value = -10986.4

$ws.Range("H132").Value = 1655.2559
$ws.Range("I132").Value = 1475.8235
$ws.Range("J132").Value = 2333.111
$ws.Range("K132").Value = 4427.470499999999
$ws.Range("L132").Value = 6999.333
$ws.Range("M132").Value = -1897.470499999999
$ws.Range("N132").Value = -12059.333
